$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New feedback entry submitted via the feedback popup - append it as the
# next row directly below the existing feedback log entries.
$ws.Range("A7").Value = "2026-02-13 19:16:41"
$ws.Range("B7").Value = "Akash"
$ws.Range("C7").Value = "Ship From Best Location"
$ws.Range("D7").Value = "Advanced Analysis"
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = "ghjghgjhg"
